$wb = $excel.ActiveWorkbook

# Helper: force a value to be stored as TEXT (even if it looks like a number),
# matching how this workbook already stores numeric-looking parameter values
# as shared strings. We briefly mark the cell as Text ("@"), write the value,
# then clear the formatting again so no residual per-cell style is left on
# the cell (keeps styles.xml effectively untouched for the cells themselves).
function Set-TextValue($ws, [string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Helper: plain numeric value.
function Set-NumValue($ws, [string]$addr, $val) {
    $ws.Range($addr).Value = $val
}

# ---- Restricciones_del_follower ----
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2
Set-TextValue $wsFollower "A2" "30.03274261603375 - x - 4.092827004219409y"
Set-TextValue $wsFollower "B2" "-28.03274261603375"
Set-TextValue $wsFollower "D2" "0.93"
Set-TextValue $wsFollower "E2" "-9.8"
Set-TextValue $wsFollower "F2" "-9.7"

# Row 3
Set-TextValue $wsFollower "A3" "1.1950000000000003 - 0.25x"
Set-TextValue $wsFollower "B3" "-3.1950000000000003"
Set-TextValue $wsFollower "D3" "0.41"
Set-TextValue $wsFollower "E3" "0"

# Row 4
Set-TextValue $wsFollower "A4" "-4.78 + x"
Set-TextValue $wsFollower "B4" "-3.2199999999999998"
Set-TextValue $wsFollower "D4" "0.7"
Set-TextValue $wsFollower "F4" "0"

# Row 5
Set-TextValue $wsFollower "A5" "-5.6 + x"
Set-TextValue $wsFollower "B5" "2.7800000000000002"
Set-TextValue $wsFollower "D5" "0.36"
Set-TextValue $wsFollower "E5" "0"
Set-TextValue $wsFollower "F5" "0"

# Row 6
Set-TextValue $wsFollower "A6" "-16.40126582278481 - 2.6582278481012658y"
Set-TextValue $wsFollower "B6" "-16.40126582278481"
Set-TextValue $wsFollower "D6" "0.79"
Set-TextValue $wsFollower "E6" "-4.0"
Set-TextValue $wsFollower "F6" "-6.3"

# ---- Punto_modificado ----
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto "A2" "4.78"
Set-TextValue $wsPunto "B2" "6.17"

# ---- Vector_bf ----
$wsBf = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $wsBf "A2" "6.90632911392405"

# ---- Vector_BF ----
$wsBF = $wb.Worksheets.Item("Vector_BF")
Set-TextValue $wsBF "A2" "-10.8"
Set-TextValue $wsBF "A3" "-51.74261603375527"

# ---- Vector_Alpha ----
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
Set-NumValue $wsAlpha "A2" 2.37
